# Apply scheduled-runner updates to Sheets/Marilith_Profits (ALC..WVR tabs)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1536
$ws.Range("I2").Value = 33.88889
$ws.Range("J2").Value = 4239.8
$ws.Range("K2").Value = 33.88889
$ws.Range("L2").Value = 4239.8
$ws.Range("M2").Value = 79.11111
$ws.Range("N2").Value = -4465.8
$ws.Range("H33").Value = 513.44446
$ws.Range("I33").Value = 379
$ws.Range("J33").Value = 984
$ws.Range("K33").Value = 379
$ws.Range("L33").Value = 984
$ws.Range("M33").Value = -150
$ws.Range("N33").Value = -1442
$ws.Range("H98").Value = 1733.7368
$ws.Range("I98").Value = 1438.4667
$ws.Range("K98").Value = 1438.4667
$ws.Range("M98").Value = 59.53330000000005
$ws.Range("H100").Value = 7999.6
$ws.Range("I100").Value = 7666.3335
$ws.Range("K100").Value = 7666.3335
$ws.Range("M100").Value = -7125.3335
$ws.Range("H107").Value = 770.7778
$ws.Range("I107").Value = 719.2857
$ws.Range("K107").Value = 719.2857
$ws.Range("M107").Value = 1200.7143
$ws.Range("H111").Value = 793.4
$ws.Range("I111").Value = 829.3333
$ws.Range("K111").Value = 2487.9999
$ws.Range("M111").Value = 579.0001000000002
$ws.Range("H112").Value = 1975.6666
$ws.Range("I112").Value = 1399.6
$ws.Range("J112").Value = 2106.5908
$ws.Range("K112").Value = 4198.799999999999
$ws.Range("L112").Value = 6319.7724
$ws.Range("M112").Value = -3090.799999999999
$ws.Range("N112").Value = -8535.7724
$ws.Range("H122").Value = 1733.7368
$ws.Range("I122").Value = 1438.4667
$ws.Range("K122").Value = 4315.4001
$ws.Range("M122").Value = -1865.4001
$ws.Range("H127").Value = 2355.7144
$ws.Range("I127").Value = 1581.6666
$ws.Range("K127").Value = 4744.9998
$ws.Range("M127").Value = 215.0002000000004
$ws.Range("H129").Value = 2297.2334
$ws.Range("I129").Value = 816.6
$ws.Range("K129").Value = 2449.8
$ws.Range("M129").Value = 2550.2
$ws.Range("H131").Value = 584.5
$ws.Range("I131").Value = 584.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1753.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 3286.5
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3412.6
$ws.Range("I132").Value = 3228.5
$ws.Range("K132").Value = 9685.5
$ws.Range("M132").Value = -7155.5
$ws.Range("H141").Value = 2572.3333
$ws.Range("I141").Value = 2572.3333
$ws.Range("K141").Value = 7716.999899999999
$ws.Range("M141").Value = -2536.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 200
$ws.Range("I18").Value = 200
$ws.Range("K18").Value = 200
$ws.Range("M18").Value = 122
$ws.Range("H45").Value = 1995.3
$ws.Range("I45").Value = 1995.3
$ws.Range("K45").Value = 1995.3
$ws.Range("M45").Value = -1618.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5325.136
$ws.Range("I134").Value = 5661.647
$ws.Range("K134").Value = 16984.941
$ws.Range("M134").Value = -14449.941
$ws.Range("H135").Value = 45999.285
$ws.Range("J135").Value = 45999.285
$ws.Range("L135").Value = 45999.285
$ws.Range("N135").Value = -56139.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1171.6285
$ws.Range("I6").Value = 1288
$ws.Range("J6").Value = 880.7
$ws.Range("K6").Value = 1288
$ws.Range("L6").Value = 880.7
$ws.Range("M6").Value = -1175
$ws.Range("N6").Value = -1106.7
$ws.Range("H23").Value = 11249.667
$ws.Range("J23").Value = 30250
$ws.Range("L23").Value = 30250
$ws.Range("N23").Value = -30730
$ws.Range("H27").Value = 11249.667
$ws.Range("J27").Value = 30250
$ws.Range("L27").Value = 30250
$ws.Range("N27").Value = -30634
$ws.Range("H58").Value = 2294.9
$ws.Range("I58").Value = 2165.9375
$ws.Range("J58").Value = 2810.75
$ws.Range("K58").Value = 2165.9375
$ws.Range("L58").Value = 2810.75
$ws.Range("M58").Value = -1962.9375
$ws.Range("N58").Value = -3216.75
$ws.Range("H132").Value = 1047.4286
$ws.Range("I132").Value = 988.6667
$ws.Range("K132").Value = 2966.0001
$ws.Range("M132").Value = -436.0001000000002
$ws.Range("H136").Value = 2294.9
$ws.Range("I136").Value = 2165.9375
$ws.Range("J136").Value = 2810.75
$ws.Range("K136").Value = 6497.8125
$ws.Range("L136").Value = 8432.25
$ws.Range("M136").Value = -3947.8125
$ws.Range("N136").Value = -13532.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 127
$ws.Range("J15").Value = 144
$ws.Range("L15").Value = 432
$ws.Range("N15").Value = -712
$ws.Range("H140").Value = 3287.3076
$ws.Range("I140").Value = 2976
$ws.Range("J140").Value = 4999.5
$ws.Range("K140").Value = 8928
$ws.Range("L140").Value = 14998.5
$ws.Range("M140").Value = -3748
$ws.Range("N140").Value = -25358.5
$ws.Range("H141").Value = 1991
$ws.Range("I141").Value = 1991
$ws.Range("K141").Value = 5973
$ws.Range("M141").Value = -793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3998.6667
$ws.Range("I126").Value = 3998.6667
$ws.Range("K126").Value = 11996.0001
$ws.Range("M126").Value = -9526.000100000001
$ws.Range("H132").Value = 1618
$ws.Range("I132").Value = 1490.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4471.9998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1941.9998
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H100").Value = 1473.25
$ws.Range("I100").Value = 964.3333
$ws.Range("K100").Value = 964.3333
$ws.Range("M100").Value = -423.3333
$ws.Range("H132").Value = 6650.2
$ws.Range("I132").Value = 7166.8887
$ws.Range("K132").Value = 21500.6661
$ws.Range("M132").Value = -18970.6661
$ws.Range("H136").Value = 3501.75
$ws.Range("I136").Value = 3501.75
$ws.Range("K136").Value = 10505.25
$ws.Range("M136").Value = -7955.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 616.3333
$ws.Range("I100").Value = 333.66666
$ws.Range("K100").Value = 667.33332
$ws.Range("M100").Value = -126.33332
$ws.Range("H107").Value = 1017.7273
$ws.Range("I107").Value = 1073.6666
$ws.Range("J107").Value = 950.6
$ws.Range("K107").Value = 3220.9998
$ws.Range("L107").Value = 2851.8
$ws.Range("M107").Value = -1300.9998
$ws.Range("N107").Value = -6691.8
